{"js": "// Eliminar capturas duplicadas: se conserva la primera imagen (la\n// \"destacada\", junto al t\u00edtulo) y se eliminan el resto de im\u00e1genes\n// repetidas (las que aparecen tras \"Introducci\u00f3n\" y \"Presentaci\u00f3n del\n// proyecto\"), junto con el p\u00e1rrafo completo que las contiene.\nconst inlinePictures = context.document.body.inlinePictures;\ninlinePictures.load(\"items\");\nawait context.sync();\n\nfor (let i = inlinePictures.items.length - 1; i >= 1; i--) {\n  inlinePictures.items[i].paragraph.delete();\n}\nawait context.sync();\n", "ps1": "# Eliminar capturas duplicadas: se conserva la primera imagen (la\n# \"destacada\", junto al t\u00edtulo) y se eliminan el resto de im\u00e1genes\n# repetidas (las que aparecen tras \"Introducci\u00f3n\" y \"Presentaci\u00f3n del\n# proyecto\"), junto con el p\u00e1rrafo completo que las contiene.\n$d = $word.ActiveDocument\n\nfor ($i = $d.InlineShapes.Count; $i -ge 2; $i--) {\n    $shape = $d.InlineShapes.Item($i)\n    $paragraph = $shape.Range.Paragraphs.Item(1)\n    $paragraph.Range.Delete()\n}\n"}
